$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest rows (old rows 10 and 11) - the sheet now only
# spans down to row 9.
$ws.Rows("10:11").Delete()

# Refresh the remaining 8 data rows (rows 2-9) with the new export of
# the support-ticket log. Every row's timestamp/content was refreshed
# and a couple of rows gained/lost answer text.

# Row 2
$ws.Cells.Item(2, 1).Value = 45541.91913194444
$ws.Cells.Item(2, 2).Value = "Нет моего вопроса"
$ws.Cells.Item(2, 3).Value = 1006569664
$ws.Cells.Item(2, 4).Value = "Roman"
$ws.Cells.Item(2, 5).Value = "Chiper"
$ws.Cells.Item(2, 6).Value = "RomanKiper"
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = ""
$ws.Cells.Item(2, 9).Value = ""
$ws.Cells.Item(2, 10).Value = ""
$ws.Cells.Item(2, 11).Value = "адлэжыдлаэыждв"

# Row 3
$ws.Cells.Item(3, 1).Value = 45541.91890046297
$ws.Cells.Item(3, 2).Value = "Помощь с подбором курса"
$ws.Cells.Item(3, 3).Value = 1006569664
$ws.Cells.Item(3, 4).Value = "Roman"
$ws.Cells.Item(3, 5).Value = "Chiper"
$ws.Cells.Item(3, 6).Value = "RomanKiper"
$ws.Cells.Item(3, 7).Value = "Отьллл"
$ws.Cells.Item(3, 8).Value = "Ооллддд"
$ws.Cells.Item(3, 9).Value = "Оолдддд"
$ws.Cells.Item(3, 10).Value = ""
$ws.Cells.Item(3, 11).Value = ""

# Row 4
$ws.Cells.Item(4, 1).Value = 45541.91875
$ws.Cells.Item(4, 2).Value = "Помощь с подбором курса"
$ws.Cells.Item(4, 3).Value = 1006569664
$ws.Cells.Item(4, 4).Value = "Roman"
$ws.Cells.Item(4, 5).Value = "Chiper"
$ws.Cells.Item(4, 6).Value = "RomanKiper"
$ws.Cells.Item(4, 7).Value = "Олллдд"
$ws.Cells.Item(4, 8).Value = "Ролллдд"
$ws.Cells.Item(4, 9).Value = "Проллл"
$ws.Cells.Item(4, 10).Value = ""
$ws.Cells.Item(4, 11).Value = ""

# Row 5
$ws.Cells.Item(5, 1).Value = 45541.91858796297
$ws.Cells.Item(5, 2).Value = "Не пришел код"
$ws.Cells.Item(5, 3).Value = 1006569664
$ws.Cells.Item(5, 4).Value = "Roman"
$ws.Cells.Item(5, 5).Value = "Chiper"
$ws.Cells.Item(5, 6).Value = "RomanKiper"
$ws.Cells.Item(5, 7).Value = ""
$ws.Cells.Item(5, 8).Value = ""
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = "Оллддд"
$ws.Cells.Item(5, 11).Value = ""

# Row 6
$ws.Cells.Item(6, 1).Value = 45541.91846064815
$ws.Cells.Item(6, 2).Value = "Не пришел код"
$ws.Cells.Item(6, 3).Value = 1006569664
$ws.Cells.Item(6, 4).Value = "Roman"
$ws.Cells.Item(6, 5).Value = "Chiper"
$ws.Cells.Item(6, 6).Value = "RomanKiper"
$ws.Cells.Item(6, 7).Value = ""
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = ""
$ws.Cells.Item(6, 10).Value = "Ооллдддд"
$ws.Cells.Item(6, 11).Value = ""

# Row 7
$ws.Cells.Item(7, 1).Value = 45541.91824074074
$ws.Cells.Item(7, 2).Value = "Не пришел код"
$ws.Cells.Item(7, 3).Value = 1006569664
$ws.Cells.Item(7, 4).Value = "Roman"
$ws.Cells.Item(7, 5).Value = "Chiper"
$ws.Cells.Item(7, 6).Value = "RomanKiper"
$ws.Cells.Item(7, 7).Value = ""
$ws.Cells.Item(7, 8).Value = ""
$ws.Cells.Item(7, 9).Value = ""
$ws.Cells.Item(7, 10).Value = "Ооолллл"
$ws.Cells.Item(7, 11).Value = ""

# Row 8
$ws.Cells.Item(8, 1).Value = 45541.9178587963
$ws.Cells.Item(8, 2).Value = "Нет моего вопроса"
$ws.Cells.Item(8, 3).Value = 1006569664
$ws.Cells.Item(8, 4).Value = "Roman"
$ws.Cells.Item(8, 5).Value = "Chiper"
$ws.Cells.Item(8, 6).Value = "RomanKiper"
$ws.Cells.Item(8, 7).Value = ""
$ws.Cells.Item(8, 8).Value = ""
$ws.Cells.Item(8, 9).Value = ""
$ws.Cells.Item(8, 10).Value = ""
$ws.Cells.Item(8, 11).Value = "Роооовлвлалал"

# Row 9
$ws.Cells.Item(9, 1).Value = 45541.9177199074
$ws.Cells.Item(9, 2).Value = "Не работает код"
$ws.Cells.Item(9, 3).Value = 1006569664
$ws.Cells.Item(9, 4).Value = "Roman"
$ws.Cells.Item(9, 5).Value = "Chiper"
$ws.Cells.Item(9, 6).Value = "RomanKiper"
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 8).Value = ""
$ws.Cells.Item(9, 9).Value = ""
$ws.Cells.Item(9, 10).Value = ""
$ws.Cells.Item(9, 11).Value = ""
